$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.132262229919434
$ws.Range("B1").Value = 3.190351009368896
$ws.Range("C1").Value = 2.608749628067017
$ws.Range("D1").Value = 2.398351907730103
$ws.Range("E1").Value = 1.915135145187378
